# Update stock prices (column C) for 2023-10 on both sheets.
$wb = $excel.ActiveWorkbook

# Sheet "个人持仓" (Personal Holdings) - rows 2..24
$ws1 = $wb.Worksheets.Item("个人持仓")
$ws1.Range("C2").Value = 47.68
$ws1.Range("C3").Value = 31.77
$ws1.Range("C4").Value = 47.58
$ws1.Range("C5").Value = 52.09
$ws1.Range("C6").Value = 30.67
$ws1.Range("C7").Value = 42.1
$ws1.Range("C8").Value = 28.54
$ws1.Range("C9").Value = 26.18
$ws1.Range("C10").Value = 27.52
$ws1.Range("C11").Value = 119.12
$ws1.Range("C12").Value = 150.41
$ws1.Range("C13").Value = 203.84
$ws1.Range("C14").Value = 0.74
$ws1.Range("C15").Value = 10.11
$ws1.Range("C16").Value = 26.85
$ws1.Range("C17").Value = 26.1
$ws1.Range("C18").Value = 15.46
$ws1.Range("C19").Value = 33.58
$ws1.Range("C20").Value = 34.79
$ws1.Range("C21").Value = 25.12
$ws1.Range("C22").Value = 95.93000000000001
$ws1.Range("C23").Value = 3.705
$ws1.Range("C24").Value = 4.48

# Sheet "家庭持仓" (Family Holdings) - rows 2..27
$ws2 = $wb.Worksheets.Item("家庭持仓")
$ws2.Range("C2").Value = 47.68
$ws2.Range("C3").Value = 31.77
$ws2.Range("C4").Value = 47.58
$ws2.Range("C5").Value = 52.09
$ws2.Range("C6").Value = 30.67
$ws2.Range("C7").Value = 42.1
$ws2.Range("C8").Value = 28.54
$ws2.Range("C9").Value = 26.18
$ws2.Range("C10").Value = 27.52
$ws2.Range("C11").Value = 119.12
$ws2.Range("C12").Value = 203.84
$ws2.Range("C13").Value = 150.41
$ws2.Range("C14").Value = 229.4
$ws2.Range("C15").Value = 0.74
$ws2.Range("C16").Value = 36.48
$ws2.Range("C17").Value = 10.11
$ws2.Range("C18").Value = 26.85
$ws2.Range("C19").Value = 26.1
$ws2.Range("C20").Value = 15.46
$ws2.Range("C21").Value = 33.58
$ws2.Range("C22").Value = 34.79
$ws2.Range("C23").Value = 25.12
$ws2.Range("C24").Value = 95.93000000000001
$ws2.Range("C25").Value = 3.705
$ws2.Range("C26").Value = 1.066
$ws2.Range("C27").Value = 4.48
